$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2025-05-07 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-05-08 Thursday", 2) | Out-Null

# Update each table cell value (table has duplicate text values, so address by row/col)
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "60÷4=15, 0"

$cell = $t.Cell(1, 2)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "29÷3=9, 2"

$cell = $t.Cell(1, 3)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "76÷8=9, 4"

$cell = $t.Cell(1, 4)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "84÷3=28, 0"

$cell = $t.Cell(1, 5)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "64÷2=32, 0"

$cell = $t.Cell(5, 1)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "24÷8=3, 0"

$cell = $t.Cell(5, 2)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "10÷8=1, 2"

$cell = $t.Cell(5, 3)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "76÷3=25, 1"

$cell = $t.Cell(5, 4)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "93÷5=18, 3"

$cell = $t.Cell(5, 5)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "37÷7=5, 2"

$cell = $t.Cell(9, 1)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "30÷9=3, 3"

$cell = $t.Cell(9, 2)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "46÷5=9, 1"

$cell = $t.Cell(9, 3)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "46÷3=15, 1"

$cell = $t.Cell(9, 4)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "50÷9=5, 5"

$cell = $t.Cell(9, 5)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "76÷9=8, 4"

$cell = $t.Cell(13, 1)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "61÷5=12, 1"

$cell = $t.Cell(13, 2)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "14÷2=7, 0"

$cell = $t.Cell(13, 3)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "53÷8=6, 5"

$cell = $t.Cell(13, 4)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "31÷3=10, 1"

$cell = $t.Cell(13, 5)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "39÷9=4, 3"

$cell = $t.Cell(17, 1)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "13÷2=6, 1"

$cell = $t.Cell(17, 2)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "28÷9=3, 1"

$cell = $t.Cell(17, 3)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "23÷4=5, 3"

$cell = $t.Cell(17, 4)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "23÷5=4, 3"

$cell = $t.Cell(17, 5)
$r = $cell.Range
$r.MoveEnd(1, -2) | Out-Null
$r.Text = "59÷8=7, 3"
